$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: SouthKorea
$ws.Range("A2").Value = "SouthKorea"
$ws.Range("B2").Value = 43942
$ws.Range("C2").Value = 0.02218477955630441
$ws.Range("D2").Value = -0.0002356190647932951
$ws.Range("E2").Value = 0.000009631854279609444
$ws.Range("F2").Value = -0.0002452509190729032
$ws.Range("G2").Value = 0.03778934979763508
$ws.Range("H2").Value = 0.9622106502023649

# Row 3: China
$ws.Range("A3").Value = "China"
$ws.Range("B3").Value = 43872
$ws.Range("C3").Value = 0.02290025071633238
$ws.Range("D3").Value = -0.000951090224821264
$ws.Range("E3").Value = -0.001865905070842499
$ws.Range("F3").Value = 0.0009148148460212372
$ws.Range("G3").Value = 0.6710151063854642
$ws.Range("H3").Value = 0.3289848936145358

# Row 4: France
$ws.Range("A4").Value = "France"
$ws.Range("B4").Value = 43914
$ws.Range("C4").Value = 0.03983587515221891
$ws.Range("D4").Value = -0.0178867146607078
$ws.Range("E4").Value = -0.02067143833938045
$ws.Range("F4").Value = 0.002784723678672646
$ws.Range("G4").Value = 0.8812796536564944
$ws.Range("H4").Value = 0.1187203463435056

# Row 5: USA
$ws.Range("A5").Value = "USA"
$ws.Range("B5").Value = 43940
$ws.Range("C5").Value = 0.06870385174884934
$ws.Range("D5").Value = -0.04675469125733822
$ws.Range("E5").Value = -0.01246165171055522
$ws.Range("F5").Value = -0.034293039546783
$ws.Range("G5").Value = 0.2665326489264187
$ws.Range("H5").Value = 0.7334673510735814

# Row 6: Spain
$ws.Range("A6").Value = "Spain"
$ws.Range("B6").Value = 43937
$ws.Range("C6").Value = 0.1050210003716739
$ws.Range("D6").Value = -0.08307183988016276
$ws.Range("E6").Value = -0.05560002378836659
$ws.Range("F6").Value = -0.02747181609179617
$ws.Range("G6").Value = 0.66930049784107
$ws.Range("H6").Value = 0.33069950215893

# Row 7: Italy
$ws.Range("A7").Value = "Italy"
$ws.Range("B7").Value = 43941
$ws.Range("C7").Value = 0.1272752828730058
$ws.Range("D7").Value = -0.1053261223814947
$ws.Range("E7").Value = -0.06796855464528848
$ws.Range("F7").Value = -0.03735756773620624
$ws.Range("G7").Value = 0.6453152656574986
$ws.Range("H7").Value = 0.3546847343425013

# Row 8: Germany - only Country and Date are populated; the rest of the
# row's figures are cleared out entirely
$ws.Range("A8").Value = "Germany"
$ws.Range("B8").Value = 43941
$ws.Range("C8:H8").ClearContents()

# Rows 9 and 10 are removed entirely from the table
$ws.Range("A9:H10").Delete()
